# "Added New Case scenario"
# Adds a new "AddNewCases" worksheet (after the existing "SearchCases" sheet),
# populates it with a header row + one data row, and makes it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after the current (only) sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "AddNewCases"

# Keep the same look & feel (Arial 10) as the rest of the workbook.
$ws2.Range("A1:E2").Font.Name = "Arial"
$ws2.Range("A1:E2").Font.Size = 10

# Header row.
$ws2.Range("A1").Value = "Cpimsid"
$ws2.Range("B1").Value = "FirstName"
$ws2.Range("C1").Value = "LastName"
$ws2.Range("D1").Value = "AssesmentDueDate"
$ws2.Range("E1").Value = "age"

# Data row.
$ws2.Range("A2").Value = "Id_cpims_1"
$ws2.Range("B2").Value = "Janani"
$ws2.Range("C2").Value = "Panchalingam"

# Store the due date as literal text "10-10-2024" (not a date serial).
$ws2.Range("D1:D2").NumberFormat = "@"
$ws2.Range("D2").Value = "10-10-2024"
$ws2.Range("D1:D2").NumberFormat = "General"

$ws2.Range("E2").Value = 33

# Make the new sheet the active / visible tab.
$ws2.Activate()
